# Deliverable1PartC.xlsx - add a new "1.3 V Controls" section (rows 20-22)
# Data is entered column-by-column (all of column A, then B, then C, then D)
# to match the order new strings were appended to the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A
$ws.Range("A20").Value = "1.3 V Controls"
$ws.Range("A21").Value = "IAW-001"
$ws.Range("A22").Value = "AD – 007"

# Column B
$ws.Range("B21").Value = "Designated admin workstation "
$ws.Range("B22").Value = "Review of PowerShell change logs"

# Column C
$ws.Range("C21").Value = "Preventive Technical preventive control "
$ws.Range("C22").Value = ": Detective Administrative control"

# Column D (wrapped, same style as the other "Brief" cells above it)
$ws.Range("D21").Value = "Create workstations only available to administrators to perform administrative actions, and prevent other workstations from completing "
$ws.Range("D22").Value = "SOC staff reviews PowerShell change logs to look for suspicious activity "

$ws.Range("D21").WrapText = $true
$ws.Range("D22").WrapText = $true

$ws.Rows.Item(21).RowHeight = 43.5
$ws.Rows.Item(22).RowHeight = 29

# Reflect the author's final cursor position/selection in the saved view
# (scrolled so row 7 is at the top, with D22 selected).
$ws.Range("D22").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
